$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two "Stage" labels in the schedule
$ws.Range("B2").Value = "Stage 1: Yes or No app which uses Face Tracking as Input"
$ws.Range("B6").Value = "Stage 3: Adding Sound & Input Sensitivity Settings"

# Minor alignment tweaks to the "Audit" column
$ws.Range("C1").HorizontalAlignment = 1
$ws.Range("C2").VerticalAlignment = -4108
$ws.Range("C5").VerticalAlignment = -4108
$ws.Range("C11").VerticalAlignment = -4108

# Move the last selected cell as recorded by Excel
$ws.Range("B14").Select() | Out-Null
